# Daily attendance processing - 2025-10-29 14:22:14
# Reorders the comma-separated list of "Recorded By" values in column G
# for the affected rows (cosmetic normalization of the recorder list order).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "backup@backdoor.com, System, system"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    8   = "backup@backdoor.com, System"
    11  = "System, dnasr281@gmail.com"
    17  = "System, dnasr281@gmail.com"
    29  = "backup@backdoor.com, System, system"
    31  = "backup@backdoor.com, System"
    32  = "backup@backdoor.com, System"
    35  = "backup@backdoor.com, System"
    38  = "System, dnasr281@gmail.com"
    44  = "System, dnasr281@gmail.com"
    56  = "backup@backdoor.com, System, system"
    58  = "backup@backdoor.com, System"
    59  = "backup@backdoor.com, System"
    62  = "backup@backdoor.com, System"
    65  = "System, dnasr281@gmail.com"
    71  = "System, dnasr281@gmail.com"
    83  = "backup@backdoor.com, System"
    84  = "backup@backdoor.com, System"
    85  = "backup@backdoor.com, System"
    90  = "admin@admin.com, dnasr281@gmail.com"
    96  = "System, dnasr281@gmail.com"
    97  = "System, dnasr281@gmail.com"
    99  = "System, dnasr281@gmail.com"
    109 = "backup@backdoor.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    116 = "admin@admin.com, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    123 = "System, dnasr281@gmail.com"
    125 = "System, dnasr281@gmail.com"
    135 = "backup@backdoor.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    142 = "admin@admin.com, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
    149 = "System, dnasr281@gmail.com"
    151 = "System, dnasr281@gmail.com"
}

foreach ($rowNum in $updates.Keys) {
    $ws.Range("G$rowNum").Value = $updates[$rowNum]
}
